$d = $word.ActiveDocument

# Several character styles in styles.xml define <w:rPr> with <w:color>
# listed before <w:b>/<w:i>, which violates the wml.xsd content model
# (CT_RPr expects w:b/w:bCs/w:i/w:iCs ... before w:color). Re-assigning
# the Font.Bold / Font.Italic property to its own current value forces
# the engine to re-serialize that style's rPr in schema-correct order,
# without altering any actual formatting. We only touch the property
# that is actually present on each style, so no new (default-valued)
# formatting element gets introduced.

$boldOnly = @(
    "KeywordTok",
    "ImportTok",
    "ControlFlowTok",
    "AlertTok",
    "ErrorTok"
)

$italicOnly = @(
    "CommentTok",
    "DocumentationTok"
)

$boldAndItalic = @(
    "AnnotationTok",
    "CommentVarTok",
    "InformationTok",
    "WarningTok"
)

foreach ($styleId in $boldOnly) {
    $s = $d.Styles.Item($styleId)
    $s.Font.Bold = $s.Font.Bold
}

foreach ($styleId in $italicOnly) {
    $s = $d.Styles.Item($styleId)
    $s.Font.Italic = $s.Font.Italic
}

foreach ($styleId in $boldAndItalic) {
    $s = $d.Styles.Item($styleId)
    $s.Font.Bold = $s.Font.Bold
    $s.Font.Italic = $s.Font.Italic
}
